# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder country names (shared-string positions shuffled upstream) ---
# Swap "Nueva Caledonia" (row 196) <-> "Belice" (row 197)
$ws.Range("A196").Value = "Belice"
$ws.Range("A197").Value = "Nueva Caledonia"

# Rotate "Seychelles" (row 209), "Groenlandia" (row 210), "Montserrat" (row 211)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Montserrat"
$ws.Range("A211").Value = "Seychelles"

# Swap "San Bartolome" (row 215) <-> "Bonaire, San Eustaquio y Saba" (row 216)
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A216").Value = "San Bartolome"

# --- Update the "last updated" timestamp footer ---
$ws.Range("A1").Value = "Datos actualizados a 19 de Mayo de 2020 a las 20:35"

# --- Update Estados Unidos (row 4) stats ---
$ws.Range("B4").Value = 1559500
$ws.Range("C4").Value = 9206
$ws.Range("E4").Value = 1106854
$ws.Range("G4").Value = 607
$ws.Range("H4").Value = 92588

# --- Update India (row 14) stats ---
$ws.Range("D14").Value = 42307
$ws.Range("E14").Value = 60838

# --- Update Israel (row 39) stats ---
$ws.Range("B39").Value = 16659
$ws.Range("C39").Value = 16
$ws.Range("D39").Value = 13435
$ws.Range("E39").Value = 2946
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 278

# --- Update stats that travelled with the name swaps ---
# Row 196 now "Belice"
$ws.Range("D196").Value = 16
$ws.Range("H196").Value = 2
# Row 197 now "Nueva Caledonia"
$ws.Range("D197").Value = 18
$ws.Range("H197").Value = 0

# Row 210 now "Montserrat"
$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1
# Row 211 now "Seychelles"
$ws.Range("D211").Value = 11
$ws.Range("H211").Value = 0
